$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 16 (shifts existing rows 16-83 down to 20-87)
$ws.Rows("16:19").Insert()

$ws.Range("A16").Value = "/DPTS"
$ws.Range("B16").Value = "List of Departments of Org"

$ws.Range("A17").Value = "/DPTS/F"
$ws.Range("B17").Value = "Find a Department with Name or ID"

$ws.Range("A18").Value = "/DPTS/MOD"
$ws.Range("B18").Value = "Modify a Department"

$ws.Range("A19").Value = "/DPTS/NEW"
$ws.Range("B19").Value = "Create a Department"

# Re-sort the whole data range (A2:B87) alphabetically by column A,
# matching a manual "sort after insert" pass.
$sortRange = $ws.Range("A2:B87")
$key1 = $ws.Range("A1")
$sortRange.Sort($key1, 1, $null, $null, 1, $null, 1, 0)

$ws.Range("A9").Select()
